$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the existing "Balcon du Jura" row now becomes index 1 (was 2)
$ws.Range("A2").Value = 1

# Row 3: Mont-Aubert
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 7040000000
$ws.Range("C3").Value = 7040
$ws.Range("D3").Value = 7000
$ws.Range("E3").Value = "Mont-Aubert"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "P"

# Row 4: Montagny - Champvent
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 7050000000
$ws.Range("C4").Value = 7050
$ws.Range("D4").Value = 7000
$ws.Range("E4").Value = "Montagny – Champvent"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "P"

# Row 5: Paquier - Donneloye
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 7060000000
$ws.Range("C5").Value = 7060
$ws.Range("D5").Value = 7000
$ws.Range("E5").Value = "Pâquier – Donneloye"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "P"

# Row 6: Yvonand
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 7100000000
$ws.Range("C6").Value = 7100
$ws.Range("D6").Value = 7000
$ws.Range("E6").Value = "Yvonand"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "P"

# Update the active selection to G7
$ws.Range("G7").Select()
